$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: item 1.0 - quantity upto date ---
$ws.Range("C8").Value = 16

# --- Row 9: Short point ---
$ws.Range("C9").Value = 51
$ws.Range("G9").Value = "'13056.00"

# --- Row 10: Medium point ---
$ws.Range("C10").Value = 79
$ws.Range("G10").Value = "'37288.00"

# --- Row 11: Long point ---
$ws.Range("C11").Value = 58
$ws.Range("G11").Value = "'38396.00"

# --- Row 12: item 2.0 - quantity upto date ---
$ws.Range("C12").Value = 7

# --- Row 13: On board ---
$ws.Range("C13").Value = 41
$ws.Range("G13").Value = "'5576.00"

# --- Row 14: P&F switch item ---
$ws.Range("C14").Value = 77
$ws.Range("G14").Value = "'1771.00"

# --- Row 15: Total ---
$ws.Range("C15").Value = 76

# --- Row 16: Add Tender Premium ---
$ws.Range("C16").Value = 6

# --- Row 17: Grand Total ---
$ws.Range("C17").Value = 65

# --- Row 19: Grand Total Rs. ---
$ws.Range("G19").Value = "'96087.00"
$ws.Range("H19").Value = "'96087.00"

# --- Row 21: NET PAYABLE AMOUNT Rs. ---
$ws.Range("G21").Value = "'96087.00"
$ws.Range("H21").Value = "'96087.00"
